$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Find.Execute("2023-04-01 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-02 Sunday", 2) | Out-Null

# Update each answer cell in the table by row/column position
# (direct Range.Text assignment avoids issues with duplicate values)
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "65×78=5070"  # was: 13×68=884
$t.Rows.Item(1).Cells.Item(2).Range.Text = "21×41=861"  # was: 99×94=9306
$t.Rows.Item(1).Cells.Item(3).Range.Text = "74×49=3626"  # was: 89×99=8811
$t.Rows.Item(1).Cells.Item(4).Range.Text = "91×29=2639"  # was: 78×96=7488
$t.Rows.Item(1).Cells.Item(5).Range.Text = "28×50=1400"  # was: 32×57=1824
$t.Rows.Item(2).Cells.Item(1).Range.Text = "90×35=3150"  # was: 60×87=5220
$t.Rows.Item(2).Cells.Item(2).Range.Text = "23×43=989"  # was: 45×97=4365
$t.Rows.Item(2).Cells.Item(3).Range.Text = "13×85=1105"  # was: 40×87=3480
$t.Rows.Item(2).Cells.Item(4).Range.Text = "90×70=6300"  # was: 63×59=3717
$t.Rows.Item(2).Cells.Item(5).Range.Text = "60×36=2160"  # was: 16×87=1392
$t.Rows.Item(3).Cells.Item(1).Range.Text = "100×22=2200"  # was: 18×73=1314
$t.Rows.Item(3).Cells.Item(2).Range.Text = "49×76=3724"  # was: 90×59=5310
$t.Rows.Item(3).Cells.Item(3).Range.Text = "89×96=8544"  # was: 19×87=1653
$t.Rows.Item(3).Cells.Item(4).Range.Text = "27×25=675"  # was: 37×46=1702
$t.Rows.Item(3).Cells.Item(5).Range.Text = "29×52=1508"  # was: 26×41=1066
$t.Rows.Item(4).Cells.Item(1).Range.Text = "83×29=2407"  # was: 49×32=1568
$t.Rows.Item(4).Cells.Item(2).Range.Text = "85×39=3315"  # was: 47×99=4653
$t.Rows.Item(4).Cells.Item(3).Range.Text = "86×24=2064"  # was: 69×91=6279
$t.Rows.Item(4).Cells.Item(4).Range.Text = "79×92=7268"  # was: 85×91=7735
$t.Rows.Item(4).Cells.Item(5).Range.Text = "49×60=2940"  # was: 23×22=506
$t.Rows.Item(5).Cells.Item(1).Range.Text = "72×11=792"  # was: 92×48=4416
$t.Rows.Item(5).Cells.Item(2).Range.Text = "40×94=3760"  # was: 87×44=3828
$t.Rows.Item(5).Cells.Item(3).Range.Text = "70×46=3220"  # was: 13×98=1274
$t.Rows.Item(5).Cells.Item(4).Range.Text = "40×52=2080"  # was: 78×67=5226
$t.Rows.Item(5).Cells.Item(5).Range.Text = "15×90=1350"  # was: 88×98=8624
$t.Rows.Item(6).Cells.Item(1).Range.Text = "13×73=949"  # was: 56×41=2296
$t.Rows.Item(6).Cells.Item(2).Range.Text = "98×27=2646"  # was: 24×72=1728
$t.Rows.Item(6).Cells.Item(3).Range.Text = "66×88=5808"  # was: 27×89=2403
$t.Rows.Item(6).Cells.Item(4).Range.Text = "84×25=2100"  # was: 20×68=1360
$t.Rows.Item(6).Cells.Item(5).Range.Text = "15×92=1380"  # was: 16×11=176
$t.Rows.Item(7).Cells.Item(1).Range.Text = "43×29=1247"  # was: 99×20=1980
$t.Rows.Item(7).Cells.Item(2).Range.Text = "100×79=7900"  # was: 93×100=9300
$t.Rows.Item(7).Cells.Item(3).Range.Text = "59×59=3481"  # was: 24×43=1032
$t.Rows.Item(7).Cells.Item(4).Range.Text = "81×90=7290"  # was: 49×58=2842
$t.Rows.Item(7).Cells.Item(5).Range.Text = "43×80=3440"  # was: 78×57=4446
$t.Rows.Item(8).Cells.Item(1).Range.Text = "88×32=2816"  # was: 17×92=1564
$t.Rows.Item(8).Cells.Item(2).Range.Text = "59×34=2006"  # was: 41×49=2009
$t.Rows.Item(8).Cells.Item(3).Range.Text = "17×22=374"  # was: 64×14=896
$t.Rows.Item(8).Cells.Item(4).Range.Text = "62×47=2914"  # was: 84×75=6300
$t.Rows.Item(8).Cells.Item(5).Range.Text = "17×23=391"  # was: 94×27=2538
$t.Rows.Item(9).Cells.Item(1).Range.Text = "15×34=510"  # was: 45×75=3375
$t.Rows.Item(9).Cells.Item(2).Range.Text = "96×22=2112"  # was: 17×27=459
$t.Rows.Item(9).Cells.Item(3).Range.Text = "17×95=1615"  # was: 76×34=2584
$t.Rows.Item(9).Cells.Item(4).Range.Text = "32×98=3136"  # was: 59×39=2301
$t.Rows.Item(9).Cells.Item(5).Range.Text = "33×52=1716"  # was: 65×83=5395
$t.Rows.Item(10).Cells.Item(1).Range.Text = "64×32=2048"  # was: 24×52=1248
$t.Rows.Item(10).Cells.Item(2).Range.Text = "82×20=1640"  # was: 16×79=1264
$t.Rows.Item(10).Cells.Item(3).Range.Text = "98×23=2254"  # was: 64×64=4096
$t.Rows.Item(10).Cells.Item(4).Range.Text = "19×12=228"  # was: 20×95=1900
$t.Rows.Item(10).Cells.Item(5).Range.Text = "33×88=2904"  # was: 71×19=1349
$t.Rows.Item(11).Cells.Item(1).Range.Text = "14×21=294"  # was: 26×10=260
$t.Rows.Item(11).Cells.Item(2).Range.Text = "40×12=480"  # was: 18×99=1782
$t.Rows.Item(11).Cells.Item(3).Range.Text = "21×33=693"  # was: 72×88=6336
$t.Rows.Item(11).Cells.Item(4).Range.Text = "10×86=860"  # was: 88×46=4048
$t.Rows.Item(11).Cells.Item(5).Range.Text = "35×87=3045"  # was: 61×34=2074
$t.Rows.Item(12).Cells.Item(1).Range.Text = "69×70=4830"  # was: 74×62=4588
$t.Rows.Item(12).Cells.Item(2).Range.Text = "50×25=1250"  # was: 74×75=5550
$t.Rows.Item(12).Cells.Item(3).Range.Text = "25×25=625"  # was: 57×62=3534
$t.Rows.Item(12).Cells.Item(4).Range.Text = "37×28=1036"  # was: 41×83=3403
$t.Rows.Item(12).Cells.Item(5).Range.Text = "73×87=6351"  # was: 42×73=3066
$t.Rows.Item(13).Cells.Item(1).Range.Text = "29×96=2784"  # was: 29×69=2001
$t.Rows.Item(13).Cells.Item(2).Range.Text = "25×49=1225"  # was: 35×93=3255
$t.Rows.Item(13).Cells.Item(3).Range.Text = "95×32=3040"  # was: 93×20=1860
$t.Rows.Item(13).Cells.Item(4).Range.Text = "49×60=2940"  # was: 33×11=363
$t.Rows.Item(13).Cells.Item(5).Range.Text = "71×82=5822"  # was: 99×74=7326
$t.Rows.Item(14).Cells.Item(1).Range.Text = "71×64=4544"  # was: 25×38=950
$t.Rows.Item(14).Cells.Item(2).Range.Text = "42×56=2352"  # was: 31×37=1147
$t.Rows.Item(14).Cells.Item(3).Range.Text = "60×48=2880"  # was: 13×23=299
$t.Rows.Item(14).Cells.Item(4).Range.Text = "74×94=6956"  # was: 11×100=1100
$t.Rows.Item(14).Cells.Item(5).Range.Text = "52×65=3380"  # was: 91×20=1820
$t.Rows.Item(15).Cells.Item(1).Range.Text = "90×85=7650"  # was: 75×74=5550
$t.Rows.Item(15).Cells.Item(2).Range.Text = "90×99=8910"  # was: 81×26=2106
$t.Rows.Item(15).Cells.Item(3).Range.Text = "49×45=2205"  # was: 12×68=816
$t.Rows.Item(15).Cells.Item(4).Range.Text = "10×12=120"  # was: 43×54=2322
$t.Rows.Item(15).Cells.Item(5).Range.Text = "88×97=8536"  # was: 28×51=1428
$t.Rows.Item(16).Cells.Item(1).Range.Text = "24×92=2208"  # was: 59×13=767
$t.Rows.Item(16).Cells.Item(2).Range.Text = "53×45=2385"  # was: 28×26=728
$t.Rows.Item(16).Cells.Item(3).Range.Text = "34×64=2176"  # was: 74×80=5920
$t.Rows.Item(16).Cells.Item(4).Range.Text = "52×94=4888"  # was: 97×55=5335
$t.Rows.Item(16).Cells.Item(5).Range.Text = "64×48=3072"  # was: 63×40=2520
$t.Rows.Item(17).Cells.Item(1).Range.Text = "32×34=1088"  # was: 12×12=144
$t.Rows.Item(17).Cells.Item(2).Range.Text = "33×46=1518"  # was: 95×74=7030
$t.Rows.Item(17).Cells.Item(3).Range.Text = "47×50=2350"  # was: 20×33=660
$t.Rows.Item(17).Cells.Item(4).Range.Text = "85×99=8415"  # was: 15×71=1065
$t.Rows.Item(17).Cells.Item(5).Range.Text = "11×12=132"  # was: 49×20=980
$t.Rows.Item(18).Cells.Item(1).Range.Text = "56×56=3136"  # was: 11×12=132
$t.Rows.Item(18).Cells.Item(2).Range.Text = "71×14=994"  # was: 86×68=5848
$t.Rows.Item(18).Cells.Item(3).Range.Text = "40×64=2560"  # was: 20×85=1700
$t.Rows.Item(18).Cells.Item(4).Range.Text = "49×49=2401"  # was: 90×20=1800
$t.Rows.Item(18).Cells.Item(5).Range.Text = "21×96=2016"  # was: 88×48=4224
$t.Rows.Item(19).Cells.Item(1).Range.Text = "24×79=1896"  # was: 11×20=220
$t.Rows.Item(19).Cells.Item(2).Range.Text = "35×89=3115"  # was: 64×64=4096
$t.Rows.Item(19).Cells.Item(3).Range.Text = "97×100=9700"  # was: 61×30=1830
$t.Rows.Item(19).Cells.Item(4).Range.Text = "33×12=396"  # was: 90×45=4050
$t.Rows.Item(19).Cells.Item(5).Range.Text = "20×11=220"  # was: 82×14=1148
$t.Rows.Item(20).Cells.Item(1).Range.Text = "77×77=5929"  # was: 34×67=2278
$t.Rows.Item(20).Cells.Item(2).Range.Text = "58×90=5220"  # was: 31×19=589
$t.Rows.Item(20).Cells.Item(3).Range.Text = "16×24=384"  # was: 56×18=1008
$t.Rows.Item(20).Cells.Item(4).Range.Text = "88×91=8008"  # was: 35×86=3010
$t.Rows.Item(20).Cells.Item(5).Range.Text = "56×17=952"  # was: 62×37=2294
